$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newNames = @{
    11 = "villageScene1"
    12 = "villageScene2"
    13 = "villageScene3"
    14 = "villageScene4"
    15 = "villageScene5"
    16 = "villageScene6"
}

foreach ($row in $newNames.Keys) {
    $name = $newNames[$row]
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $name
}

$ws.Cells.Item(9,1).EntireRow.AutoFit()
$ws.Rows.Item(10).RowHeight = 54.4

$ws.Range("D18").Select()
